# Update "想去人数" (interested-in count) values on the "展览" and "全部类型"
# sheets to reflect the regenerated data output.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F2").Value = 577
$sheetExhibition.Range("F3").Value = 192
$sheetExhibition.Range("F4").Value = 395
$sheetExhibition.Range("F5").Value = 436
$sheetExhibition.Range("F7").Value = 2457
$sheetExhibition.Range("F8").Value = 423
$sheetExhibition.Range("F9").Value = 6437
$sheetExhibition.Range("F11").Value = 415

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 577
$sheetAll.Range("F3").Value = 192
$sheetAll.Range("F4").Value = 395
$sheetAll.Range("F5").Value = 436
$sheetAll.Range("F9").Value = 2457
$sheetAll.Range("F10").Value = 423
$sheetAll.Range("F11").Value = 6437
$sheetAll.Range("F13").Value = 415
